$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Color constants (Excel VBA-style BGR-packed integers)
$green = 65280   # RGB(0,255,0)
$red   = 13311   # RGB(255,51,0)

# Highlight the "MU/P $2" marginal-utility-per-dollar cells that are equal
# (the equimarginal/optimal allocation points) in green, and the matching
# "not suitable" extra unit in red, for both goods' columns.
$ws.Range("D3").Font.Bold = $true
$ws.Range("D3").Interior.Color = $green

$ws.Range("D7").Font.Bold = $true
$ws.Range("D7").Interior.Color = $red

$ws.Range("H11").Font.Bold = $true
$ws.Range("H11").Interior.Color = $green

$ws.Range("H12").Font.Bold = $true
$ws.Range("H12").Interior.Color = $red

# Add a small legend below the table explaining the highlight colors.
$ws.Range("C15").Value = "Suitable"
$ws.Range("C15").Font.Bold = $true
$ws.Range("D15").Interior.Color = $green

$ws.Range("C16").Value = "Not Suitable"
$ws.Range("C16").Font.Bold = $true
$ws.Range("D16").Interior.Color = $red
